# Daily update at 8 AM UTC
# Appends the new day's results as row 23, and moves the "latest row"
# date-only number format down from row 22 to the newly appended row 23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 was previously the last row (shown with a date-only format);
# now that it's no longer the last row, give it the regular
# date+time number format used by all the other non-final rows.
$ws.Range("A22").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data in row 23.
$ws.Range("A23").Value = 45763
$ws.Range("B23").Value = 92
$ws.Range("C23").Value = 91
$ws.Range("D23").Value = 91

# Row 23 is now the last row, so it gets the date-only number format.
$ws.Range("A23").NumberFormat = "YYYY-MM-DD"
